$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AHB-Diff")

# Rename the header row: "_old" columns become "_FV2210", "_new" columns become "_FV2304"
$fv2210Headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)
$fv2304Headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

# Columns A-J carry the "_old" -> "_FV2210" headers
for ($i = 0; $i -lt $fv2210Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2210Headers[$i]
}

# Column K is "diff" - unchanged

# Columns L-U carry the "_new" -> "_FV2304" headers
for ($i = 0; $i -lt $fv2304Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2304Headers[$i]
}

# Freeze the header row (row 1)
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the used range into an Excel Table ("Table1") with banded rows
$dataRange = $ws.Range("A1:U56")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
